$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("M19P")
Write-Host $ws.Range("A5").Value
Write-Host $ws.Range("O5").Value
$ws.Range("O5").Value = "tf"
$ws.Range("O6").Value = 0.34615384615384615
